$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K7").Value = -0.763888888888889
$ws.Range("K12").Value = 12.41429539295394
$ws.Range("K13").Value = 13.75752314814816
$ws.Range("K15").Value = 21.79166666666666
$ws.Range("K16").Value = 13.75752314814816
$ws.Range("K23").Value = 20.68981481481483
$ws.Range("K26").Value = -0.763888888888889
$ws.Range("K27").Value = 12.67039049919483
$ws.Range("K28").Value = 19.65277777777778
